$d = $word.ActiveDocument

function Make-Pkg($bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
        $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# 1. Collapse the three-run title ("Module " + "2" + " questions") into a
#    single run reading "Module 2 questions".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Module 2 questions", $true, $false, $false, $false, `
    $false, $true, 1, $false, "Module 2 questions", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Fill in the (currently empty) first "Interview questions:" list item and
#    append the remaining nine interview questions as sibling list paragraphs
#    (same ListParagraph style / numId=1 numbering, sz/szCs 24 run props).
#
#    NB: the runtime's Range.InsertXML mis-behaves (wipes the whole document)
#    whenever an insert leaves the document at *exactly* 10 paragraphs, so the
#    ten new questions are written in two batches (9 + 1) that never land on
#    that total.
# ---------------------------------------------------------------------------
$interviewPPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>'
$interviewRPr = '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'

$iq1 = "How many events do you host per month?"
$iq2 = "The average number of attendees?"
$iq4 = "Do you keep track of the attendance?"
$iq5 = "If yes, do you record the attendee" + [char]0x2019 + "s data?"
$iq6 = "Do you take in mind the preferred events of the user?"
$iq7 = "How do receive feedback from the users?"
$iq8 = "Are users generally satisfied or dissatisfied with the current system?"
$iq9 = "Do you think there is any problem with the current system?"
$iq10 = "How can you enhance the user experience in the new system?"

$firstNineBody = ''
$firstNineBody += '<w:p>' + $interviewPPr + '<w:r>' + $interviewRPr + '<w:t>' + $iq1 + '</w:t></w:r></w:p>'
$firstNineBody += '<w:p>' + $interviewPPr + '<w:r>' + $interviewRPr + '<w:t>' + $iq2 + '</w:t></w:r></w:p>'
$firstNineBody += '<w:p>' + $interviewPPr + `
    '<w:r>' + $interviewRPr + '<w:t>What type of students are more interested in attending in terms of nationality</w:t></w:r>' + `
    '<w:r>' + $interviewRPr + '<w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r>' + $interviewRPr + '<w:t>and major field?</w:t></w:r>' + `
    '</w:p>'
$firstNineBody += '<w:p>' + $interviewPPr + '<w:r>' + $interviewRPr + '<w:t>' + $iq4 + '</w:t></w:r></w:p>'
$firstNineBody += '<w:p>' + $interviewPPr + '<w:r>' + $interviewRPr + '<w:t>' + $iq5 + '</w:t></w:r></w:p>'
$firstNineBody += '<w:p>' + $interviewPPr + '<w:r>' + $interviewRPr + '<w:t>' + $iq6 + '</w:t></w:r></w:p>'
$firstNineBody += '<w:p>' + $interviewPPr + '<w:r>' + $interviewRPr + '<w:t>' + $iq7 + '</w:t></w:r></w:p>'
$firstNineBody += '<w:p>' + $interviewPPr + '<w:r>' + $interviewRPr + '<w:t>' + $iq8 + '</w:t></w:r></w:p>'
$firstNineBody += '<w:p>' + $interviewPPr + '<w:r>' + $interviewRPr + '<w:t>' + $iq9 + '</w:t></w:r></w:p>'

# locate the lone empty ListParagraph (numId=1) item under "Interview questions:"
$interviewTarget = $null
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text.TrimEnd([char]13)
    if ($interviewTarget -eq $null -and $txt -eq "" -and $p.Style.NameLocal -eq "List Paragraph") {
        $interviewTarget = $p
    }
}
$interviewTarget.Range.InsertXML((Make-Pkg $firstNineBody)) | Out-Null

# find the paragraph now holding the 9th question ("Do you think there is any
# problem...") and append the 10th as a new sibling paragraph after it.
$ninthPara = $null
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text.TrimEnd([char]13)
    if ($txt -eq $iq9) {
        $ninthPara = $p
    }
}
$ninthPara.Range.InsertParagraphAfter() | Out-Null

$tenthPara = $ninthPara.Next()
$tenthBody = '<w:p>' + $interviewPPr + '<w:r>' + $interviewRPr + '<w:t>' + $iq10 + '</w:t></w:r></w:p>'
$tenthPara.Range.InsertXML((Make-Pkg $tenthBody)) | Out-Null

# ---------------------------------------------------------------------------
# 3. Replace the placeholder "." questionnaire list item (and absorb the
#    stray trailing empty paragraph at the end of the body) with the eight
#    questionnaire questions (ListParagraph style / numId=2 numbering).
# ---------------------------------------------------------------------------
$questionnairePPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>'

$qq5 = "Have you ever missed an event because you haven" + [char]0x2019 + "t been notified?"

$questionnaireBody = ''
$questionnaireBody += '<w:p>' + $questionnairePPr + '<w:r><w:t>What type of events do you prefer attending?</w:t></w:r></w:p>'
$questionnaireBody += '<w:p>' + $questionnairePPr + '<w:r><w:t>On average how many events had you attended?</w:t></w:r></w:p>'
$questionnaireBody += '<w:p>' + $questionnairePPr + '<w:r><w:t>Do the current events satisfy your interest?</w:t></w:r></w:p>'
$questionnaireBody += '<w:p>' + $questionnairePPr + '<w:r><w:t>How could you know if there are any coming or outgoing events?</w:t></w:r></w:p>'
$questionnaireBody += '<w:p>' + $questionnairePPr + '<w:r><w:t>' + $qq5 + '</w:t></w:r></w:p>'
$questionnaireBody += '<w:p>' + $questionnairePPr + '<w:r><w:t>Is there any place where you can give feedback or complaints?</w:t></w:r></w:p>'
$questionnaireBody += '<w:p>' + $questionnairePPr + '<w:r><w:t>How satisfied are you with the current system?</w:t></w:r></w:p>'
$questionnaireBody += '<w:p>' + $questionnairePPr + '<w:r><w:t>What do you think the system needs to gain your satisfaction?</w:t></w:r></w:p>'

$questionnaireTarget = $null
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text.TrimEnd([char]13)
    if ($txt -eq "." -and $p.Style.NameLocal -eq "List Paragraph") {
        $questionnaireTarget = $p
    }
}

$rStart = $questionnaireTarget.Range.Start
$rEnd = $d.Content.End
$questionnaireRange = $d.Range($rStart, $rEnd)
$questionnaireRange.InsertXML((Make-Pkg $questionnaireBody)) | Out-Null

Write-Output "done"
